# Simulated and logged 2021 conference championships
# Appends the new game's per-play data to the running season strings and
# bumps the season-total aggregate cells on OFF / DEF / ST / TURNS / PEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append the new game's per-play yardage log (rush/pass,
# offense/defense) to the running space-separated strings.
# ---------------------------------------------------------------------
$yds = $wb.Worksheets.Item("YDS")

$yds.Range("B2").Value = $yds.Range("B2").Value() + " 0 3 3 1 3 3 2 4 6 9 3 1 1 5 -1 0 -1 9 -1 0"
$yds.Range("B3").Value = $yds.Range("B3").Value() + " 12 18 2 31 44 11 6 13 11 13 14 11 16 21 12 -3"
$yds.Range("C2").Value = $yds.Range("C2").Value() + " 3 14 0 0 8 4 4 3 6 14 -3 1 1 9 2 5 2 -2 0 2 0 0 3 1 -1"
$yds.Range("C3").Value = $yds.Range("C3").Value() + " 3 3 17 2 13 10 5 6 11 7 15 16 13 9 13 26 9 7 11 20 9 11 29 16 7 -4 6 8 6 25 8"

# ---------------------------------------------------------------------
# OFF sheet: bump the season totals with the new game's numbers.
# ---------------------------------------------------------------------
$off = $wb.Worksheets.Item("OFF")

$off.Range("B2").Value = 6
$off.Range("C2").Value = 288
$off.Range("F2").Value = 88
$off.Range("G2").Value = 76
$off.Range("J2").Value = 42

$off.Range("C3").Value = 190
$off.Range("E3").Value = 45
$off.Range("F3").Value = 114
$off.Range("G3").Value = 31
$off.Range("H3").Value = 38
$off.Range("I3").Value = 70
$off.Range("J3").Value = 49
$off.Range("L3").Value = 342
$off.Range("M3").Value = 226
$off.Range("Q3").Value = 732

# ---------------------------------------------------------------------
# DEF sheet: bump the season totals with the new game's numbers.
# ---------------------------------------------------------------------
$def = $wb.Worksheets.Item("DEF")

$def.Range("C2").Value = 236
$def.Range("E2").Value = 17
$def.Range("F2").Value = 76
$def.Range("G2").Value = 72
$def.Range("I2").Value = 8
$def.Range("J2").Value = 36
$def.Range("N2").Value = 44
$def.Range("O2").Value = 28

$def.Range("C3").Value = 229
$def.Range("E3").Value = 52
$def.Range("F3").Value = 130
$def.Range("G3").Value = 36
$def.Range("H3").Value = 41
$def.Range("I3").Value = 81
$def.Range("J3").Value = 78
$def.Range("L3").Value = 414
$def.Range("M3").Value = 280
$def.Range("Q3").Value = 774

# ---------------------------------------------------------------------
# ST sheet: bump the special-teams season totals, and append the new
# game's per-kick logs to the running strings.
# ---------------------------------------------------------------------
$st = $wb.Worksheets.Item("ST")

$st.Range("B2").Value = 99
$st.Range("D2").Value = 69
$st.Range("F2").Value = 592
$st.Range("G2").Value = 577
$st.Range("J2").Value = 292
$st.Range("K2").Value = 277
$st.Range("B3").Value = 43

$st.Range("B4").Value = $st.Range("B4").Value() + " 66"
$st.Range("B5").Value = $st.Range("B5").Value() + " 24"
$st.Range("B6").Value = $st.Range("B6").Value() + " 15 11"
$st.Range("D3").Value = $st.Range("D3").Value() + " 44 37 31 35 41"
$st.Range("D4").Value = $st.Range("D4").Value() + " 0 0 0 0 0"
$st.Range("D5").Value = $st.Range("D5").Value() + " 9"

# ---------------------------------------------------------------------
# TURNS sheet: bump the Road turnover totals.
# ---------------------------------------------------------------------
$turns = $wb.Worksheets.Item("TURNS")

$turns.Range("B3").Value = 10
$turns.Range("C3").Value = 7
$turns.Range("E3").Value = 13

# ---------------------------------------------------------------------
# PEN sheet: bump the penalty totals.
# ---------------------------------------------------------------------
$pen = $wb.Worksheets.Item("PEN")

$pen.Range("B2").Value = 25
$pen.Range("D4").Value = 21
